# Kingdom Hearts Recoded - "Initial update for checkpoint code"
# Fill in Wonderland location rows on the "Locations" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

# Curly apostrophe (U+2019) used by "The Queen's Castle", matching the
# typographic quote already used elsewhere in the workbook.
$rsquo = [char]0x2019

# Data for rows that already had an index (A) and formula (F) defined,
# but were missing the B (areaId), C (worldId), D (name) and E (display)
# lookup values. Also includes brand-new rows 28-31 which need the same
# treatment, plus blank placeholder rows 32-36 that only need A & F.
$rows = @(
    @{Row=23; B="0x5";  C="0x2"; D="WonderlandAimless";        E="Aimless Path"},
    @{Row=24; B="0x4";  C="0x2"; D="WonderlandMaze";           E="Hedge Maze"},
    @{Row=25; B="0x7";  C="0x2"; D="WonderlandBizarreSmall";   E="Bizarre Room"},
    @{Row=26; B="0x1";  C="0x2"; D="WonderlandBizarreBig";     E="Bizarre Room"},
    @{Row=27; B="0x3";  C="0x2"; D="WonderlandTea";            E="Tea Party Garden"},
    @{Row=28; B="0x2";  C="0x2"; D="WonderlandQueen";          E="The Queen$($rsquo)s Castle"},
    @{Row=29; B="0xa";  C="0x2"; D="WonderlandKeyholeFirst";   E="Keyhole / The Falsewood"},
    @{Row=30; B="0xb";  C="0x2"; D="WonderlandKeyholeSecond";  E="Keyhole / Fleeting Forest"},
    @{Row=31; B="0xc";  C="0x2"; D="WonderlandKeyholeTerminus"; E="Keyhole / Terminus"}
)

foreach ($r in $rows) {
    $row = $r.Row
    # Rows 28-31 are brand new (previously absent), so column A (the
    # running index) and the F formula have to be (re)created for them
    # too, not just B:E.
    $ws.Cells.Item($row, 1).Value = $row - 2
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $formula = '=_xlfn.CONCAT( ,A' + $row + ',": { ""worldId"": ",C' + $row + ',", ""name"": """,D' + $row + ',""", ""display"": """,E' + $row + ',""", ""areaId"": ",B' + $row + ',", },")'
    $ws.Range("F$row").Formula = $formula
}

# New blank rows 32-36: only column A (index) and the dragged-down F
# formula are populated; B:E remain empty for now.
for ($row = 32; $row -le 36; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
    $formula = '=_xlfn.CONCAT( ,A' + $row + ',": { ""worldId"": ",C' + $row + ',", ""name"": """,D' + $row + ',""", ""display"": """,E' + $row + ',""", ""areaId"": ",B' + $row + ',", },")'
    $ws.Range("F$row").Formula = $formula
}

Write-Output "Locations sheet updated through row 36"
